$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Locate the "Make PNG transparent" paragraph in the TODO list and replace
#    it (and append the new items after it) with the full updated block:
#      Make PNG transparent
#      <blank>
#      #12
#      Make game over text
#      <blank>
#      <blank>
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Make PNG transparent*") {
        $target = $p
        break
    }
}

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newXml  = "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:t>Make PNG transparent</w:t></w:r></w:p>"
$newXml += "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr></w:p>"
$newXml += "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>#12</w:t></w:r></w:p>"
$newXml += "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>Make game over text</w:t></w:r></w:p>"
$newXml += "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr></w:pPr></w:p>"
$newXml += "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/></w:pPr></w:p>"

[void]$target.Range.InsertXML($newXml)

# ---------------------------------------------------------------------------
# 2) The cached page-number field in the header ("PAGE \* MERGEFORMAT")
#    should read 4 instead of 7.
# ---------------------------------------------------------------------------
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)
    $headers = $section.Headers
    for ($h = 1; $h -le $headers.Count; $h++) {
        $hdr = $headers.Item($h)
        if ($hdr.Exists) {
            $flds = $hdr.Range.Fields
            for ($fi = 1; $fi -le $flds.Count; $fi++) {
                $fld = $flds.Item($fi)
                if ($fld.Code.Text -like "*PAGE*" -and $fld.Result.Text -eq "7") {
                    $fld.Result.Text = "4"
                }
            }
        }
    }
}
